# send_list.xlsx — "ok to excel file"
#
# The sheet tracks candidates with an "ok" flag in column A (header row 1).
# This edit toggles which rows are flagged "ok":
#   - row 2 (Валентина/Абишев) is no longer "ok" -> clear it
#   - row 3 (Иван Петрович/Божевольнов) is now "ok" -> set it
#   - row 4 (Сергей и Елена/Степушин) stays "ok" -> untouched
#   - row 5 (Родимина) is now "ok" -> set it
# and leaves the active selection on A7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").ClearContents()
$ws.Range("A3").Value = "ok"
$ws.Range("A5").Value = "ok"

$ws.Range("A7").Select()
